$d = $word.ActiveDocument

# 1. Title "Capítulo 1. Introducción." -> "Capítulo 1. Introducción" (merge runs, drop trailing period)
$d.Content.Find.Execute("Capítulo 1. Introducción.", $true, $false, $false, $false, $false, $true, 1, $false, "Capítulo 1. Introducción", 2) | Out-Null

# 2. Paragraph "En este capítulo se " + "van a llevar..." merge runs (same text) - no text change needed,
#    but we still touch it via Find/Replace to force a run merge.
$d.Content.Find.Execute("En este capítulo se van a llevar a cabo una introducción al trabajo en la que se describirán los objetivos principales del trabajo y se detallará el proceso para desarrollarlo. Por último, se describirá la estructura de esta memoria.", $true, $false, $false, $false, $false, $true, 1, $false, "En este capítulo se van a llevar a cabo una introducción al trabajo en la que se describirán los objetivos principales del trabajo y se detallará el proceso para desarrollarlo. Por último, se describirá la estructura de esta memoria.", 2) | Out-Null

# 3. "1.1. Objetivos." -> "1.1. Objetivos" (drop trailing period)
$d.Content.Find.Execute("1.1. Objetivos.", $true, $false, $false, $false, $false, $true, 1, $false, "1.1. Objetivos", 2) | Out-Null

# 4. "1.2. Desarrollo del trabajo." -> "1.2. Desarrollo del trabajo" (drop trailing period)
$d.Content.Find.Execute("1.2. Desarrollo del trabajo.", $true, $false, $false, $false, $false, $true, 1, $false, "1.2. Desarrollo del trabajo", 2) | Out-Null

# 5. "1.3. Estructura de la memoria." -> drop trailing period and split into two runs:
#    "1.3. E" / "structura de la memoria"
$rng5 = $d.Content
$rng5.Find.Execute("1.3. Estructura de la memoria.") | Out-Null
$start5 = $rng5.Start
$end5 = $rng5.End

# remove the trailing period
$d.Range($end5 - 1, $end5).Delete() | Out-Null

# split the run after "1.3. E" (6 characters in) by inserting then removing a paragraph mark;
# this leaves the text as two adjacent runs instead of Word's usual re-merge of identical runs
$splitPoint5 = $start5 + 6
$d.Range($splitPoint5, $splitPoint5).InsertParagraphAfter() | Out-Null
$d.Range($splitPoint5, $splitPoint5 + 1).Delete() | Out-Null
